$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark (it will be re-created at the
#    new edit location, matching Word's normal behaviour of moving the
#    "last edit" bookmark to wherever text was most recently typed).
try {
    $old = $d.Bookmarks.Item("_GoBack")
    $old.Delete()
} catch {
}

# 2. Fill in the trailing empty list paragraph with the new changelog entry.
$p = $d.Paragraphs.Last
$p.Range.Text = "When de player dies, wait a short time before restarting the gameX"
$p.Range.LanguageID = "en-US"

# 3. Insert a collapsed "_GoBack" bookmark right after the new text. Using a
#    throwaway trailing character lets us place a truly collapsed bookmark
#    (rather than one that wraps the whole paragraph), then we remove the
#    placeholder character.
$paraEnd = $p.Range.End - 2
$markPos = $d.Range($paraEnd, $paraEnd)
$d.Bookmarks.Add("_GoBack", $markPos)

$placeholder = $d.Range($paraEnd, $paraEnd + 1)
$placeholder.Text = ""
